$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "REGALO" header in F1 - reuse E1's header formatting (bold, centered)
# so the new column matches the look of the existing header row.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "REGALO"

# Fill in the new "REGALO" data column (F2:F11)
$ws.Range("F2").Value  = "FLORES"
$ws.Range("F3").Value  = "GASEOSA"
$ws.Range("F4").Value  = "PERRO"
$ws.Range("F5").Value  = "ESCOBA"
$ws.Range("F6").Value  = "ROPA"
$ws.Range("F7").Value  = "ROPA"
$ws.Range("F8").Value  = "CALCETINES"
$ws.Range("F9").Value  = "ZAPATOS"
$ws.Range("F10").Value = "CALCETINES"
$ws.Range("F11").Value = "PERRO"

# Match the author's final selection in the sheet
$ws.Range("H8").Select()
